# Update user/student records so students (and user accounts) can access
# the sheet used to create attendance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "i"
$ws.Range("A3").Value = "b"
$ws.Range("A4").Value = "c"
$ws.Range("A5").Value = "d"

$ws.Range("B2").Value = "i@ee.com"
$ws.Range("B3").Value = "b@ee.com"
$ws.Range("B4").Value = "c@o.com"
$ws.Range("B5").Value = "d@e.com"

# Move the active selection to B5, matching the saved view state.
$ws.Range("B5").Select()
